# Adverse.docx template update:
#   {#reputationalIndicators} -> {#reputationalIndicatorsRich}
#   {.}                       -> {prefix}{authorLink}{suffix}
#   {/reputationalIndicators} -> {/reputationalIndicatorsRich}
#   "Fonti (Autore o testata):" / {#indicatorSources} / {/indicatorSources} -> "" (blanked out)

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    # wdFindContinue(1) + wdReplaceAll(2) over the whole document body.
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Rename the reputational-indicators merge field block to the "rich" variant.
Replace-Text "{#reputationalIndicators}" "{#reputationalIndicatorsRich}"
Replace-Text "{.}" "{prefix}{authorLink}{suffix}"
Replace-Text "{/reputationalIndicators}" "{/reputationalIndicatorsRich}"

# Clear the old "Fonti (Autore o testata)" / indicatorSources loop text that
# used to sit after the final section break.
Replace-Text "Fonti (Autore o testata):" ""
Replace-Text "{#indicatorSources}" ""
Replace-Text "{/indicatorSources}" ""
